$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the next empty row in column A (mirrors appending a new log entry)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2025-08-25 22:44:36"
